$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TEAM_AVAIL")

# --- Assign / Unassign team members ---
# Rows 2-10: mark as "leave" (Unassigned)
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 9).Value = "leave"
}

# Rows 11-30: mark as "office" (Assigned)
for ($r = 11; $r -le 30; $r++) {
    $ws.Cells.Item($r, 9).Value = "office"
}

# --- Resize columns to fit their content (best fit) ---
$ws.Columns.Item(1).ColumnWidth = 7.285714285714286
$ws.Columns.Item(2).ColumnWidth = 7.857142857142857
$ws.Columns.Item(3).ColumnWidth = 6.857142857142857
$ws.Columns.Item(4).ColumnWidth = 9.714285714285714
$ws.Columns.Item(5).ColumnWidth = 10.142857142857142
$ws.Columns.Item(6).ColumnWidth = 11.428571428571429
$ws.Columns.Item(7).ColumnWidth = 26.428571428571427
$ws.Columns.Item(8).ColumnWidth = 6.142857142857143
$ws.Columns.Item(9).ColumnWidth = 5.285714285714286

# --- Restore the active selection below the data ---
$ws.Range("D31").Select()
